# Append new rows of filtering-dot statistics to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @("hsa_filtering_dot\Biosynthesis of amino acids.dot", 146, 434, 244, 288, 98, -43.78, -65.97),
    @("hsa_filtering_dot\Biosynthesis of cofactors.dot", 326, 886, 452, 560, 126, -48.98, -77.5),
    @("hsa_filtering_dot\Carbon metabolism.dot", 181, 537, 317, 356, 136, -40.97, -61.8),
    @("hsa_filtering_dot\Nucleotide metabolism.dot", 214, 766, 468, 552, 254, -38.9, -53.99),
    @("hsa_filtering_dot\Purine.dot", 228, 772, 466, 544, 238, -39.64, -56.25),
    @("hsa_filtering_dot\test.dot", 13, 34, 16, 21, 3, -52.94, -85.70999999999999),
    @("mmu_filtering_dot\Biosynthesis of amino acids.dot", 150, 447, 253, 297, 103, -43.4, -65.31999999999999),
    @("mmu_filtering_dot\Biosynthesis of cofactors.dot", 329, 894, 456, 565, 127, -48.99, -77.52),
    @("mmu_filtering_dot\Carbon metabolism.dot", 181, 537, 317, 356, 136, -40.97, -61.8),
    @("mmu_filtering_dot\Nucleotide metabolism.dot", 214, 766, 468, 552, 254, -38.9, -53.99),
    @("mmu_filtering_dot\Purine.dot", 234, 789, 475, 555, 241, -39.8, -56.58)
)

$startRow = 7
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
